# CI: Sync Excel from SVN to Git
#
# The "道具描述总表" (Item) sheet's column D used to be a single
# "use_action" field. It is replaced by four new columns:
#   D: unlock_lv   E: rarity   F: show_max_stacking   G: price
# Every data row gets an unlock_lv of 1, and the two new "price"
# rows (5 & 6) get their price values in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- header row 2: field names for the new columns ---
$ws.Range("D2").Value = "unlock_lv"

$ws.Range("E2").Value = "rarity"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial($xlPasteFormats)

$ws.Range("F2").Value = "show_max_stacking"
$ws.Range("G2").Value = "price"

# --- header row 1 / row 3-8: blank E column cells, matching formats ---
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial($xlPasteFormats)
$ws.Range("D1").ClearContents()

$ws.Range("C4").Copy()
$ws.Range("E4").PasteSpecial($xlPasteFormats)

$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial($xlPasteFormats)
$ws.Range("E5").PasteSpecial($xlPasteFormats)
$ws.Range("E6").PasteSpecial($xlPasteFormats)
$ws.Range("E7").PasteSpecial($xlPasteFormats)
$ws.Range("E8").PasteSpecial($xlPasteFormats)

# --- data rows: unlock_lv = 1 (column D), re-formatted like column C ---
$ws.Range("D3").Value = "1"
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteFormats)

$ws.Range("D4").Value = "1"

$ws.Range("D5").Value = "1"
$ws.Range("D6").Value = "1"
$ws.Range("D7").Value = $null
$ws.Range("D8").Value = $null
$ws.Range("C3").Copy()
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("D7").PasteSpecial($xlPasteFormats)
$ws.Range("D8").PasteSpecial($xlPasteFormats)
$ws.Range("D5").Value = "1"
$ws.Range("D6").Value = "1"

# --- price column for the two new rows ---
$ws.Range("B5").Value = "1000000"
$ws.Range("B6").Value = "1000001"

$excel.CutCopyMode = $false

# --- column widths (character units - closest achievable to the saved file) ---
$ws.Range("B:B").ColumnWidth = 11.5000001
$ws.Range("C:E").ColumnWidth = 12.642857242857142
$ws.Range("F:F").ColumnWidth = 16.5000001

# --- selection matches the saved view ---
$ws.Range("D8").Select()
